# Retrain run: move the Loss/Accuracy figures from F:G into D:E (new
# artifacts only report two metrics now) and clear out the old F:G cells,
# keeping their number formatting in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Loss (D) / Accuracy (E) values per row, freshly produced by retraining.
$newValues = @{
    2  = @(0.72960686683654696, 0.72545754909515303)
    3  = @(0.38250002264976501, 0.87104827165603604)
    4  = @(0.364557415246963,   0.87437605857849099)
    5  = @(0.582988440990448,   0.78702163696288996)
    6  = @(0.30606013536453203, 0.89767056703567505)
    7  = @(0.37662041187286299, 0.87687188386917103)
    8  = @(0.53930389881134,    0.80782032012939398)
    9  = @(0.32160222530364901, 0.90432614088058405)
    10 = @(0.295376986265182,   0.91763728857040405)
    11 = @(0.49054548144340498, 0.83028286695480302)
}

foreach ($r in 2..11) {
    $vals = $newValues[$r]

    # Carry the existing F/G number formatting (centered, "0.00_ ") over to
    # D/E before writing the new values there.
    $ws.Range("F$r").Copy()
    $ws.Range("D$r").PasteSpecial(-4122)
    $ws.Range("G$r").Copy()
    $ws.Range("E$r").PasteSpecial(-4122)

    $ws.Range("D$r").Value2 = $vals[0]
    $ws.Range("E$r").Value2 = $vals[1]

    # Old columns are retired but keep their style/formatting in place.
    $ws.Range("F$r").ClearContents()
    $ws.Range("G$r").ClearContents()
}

$ws.Application.CutCopyMode = $false

# Match the retired columns' width to the new D column so the now-empty
# F:G still line up visually.
$ws.Range("F1:G1").EntireColumn.ColumnWidth = $ws.Range("D1").EntireColumn.ColumnWidth

# Final selection left on the sheet after the edit.
[void]$ws.Range("E11").Select()
